# Update the cryptos list sheet with the latest scraped values.
# All data cells in this sheet are plain text (prices/links/coin names are
# stored as strings, not numbers) - force the "Price" column to stay text so
# Excel's automatic number detection doesn't silently convert e.g. "581.14"
# into a floating point number (which would also corrupt values such as
# "0.0000250" into scientific notation).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple per-row D/E (price / volume-1h) updates -----------------------
# Map: row number -> @{ D = new price (or $null if unchanged); E = new pct (or $null if unchanged) }
$updates = @{
    2  = @{ D = "66.917.66";  E = "  +2.31%  " }
    3  = @{ D = "3.120.10";   E = "  +5.76%  " }
    4  = @{ D = $null;        E = "  +0.03%  " }
    5  = @{ D = "581.14";     E = "  +2.04%  " }
    6  = @{ D = "172.28";     E = "  +7.47%  " }
    7  = @{ D = $null;        E = "  +0.00%  " }
    8  = @{ D = "3.114.62";   E = "  +5.65%  " }
    9  = @{ D = "0.524";      E = "  +1.41%  " }
    10 = @{ D = $null;        E = "  -3.13%  " }
    11 = @{ D = $null;        E = "  +4.17%  " }
    12 = @{ D = $null;        E = "  +5.01%  " }
    13 = @{ D = "0.0000250";  E = "  +2.36%  " }
    14 = @{ D = "37.21";      E = "  +7.87%  " }
    15 = @{ D = $null;        E = "  -0.07%  " }
    16 = @{ D = "3.631.45";   E = "  +5.58%  " }
    17 = @{ D = "66.947.02";  E = "  +2.34%  " }
    18 = @{ D = "7.22";       E = "  +2.94%  " }
    19 = @{ D = "3.114.97";   E = "  +5.62%  " }
    20 = @{ D = "16.22";      E = "  +3.14%  " }
    21 = @{ D = "484.34";     E = "  +8.93%  " }
    22 = @{ D = "0.718";      E = "  +3.40%  " }
    23 = @{ D = "7.56";       E = "  +3.67%  " }
    24 = @{ D = "84.21";      E = "  +2.41%  " }
    25 = @{ D = "2.37";       E = "  +5.87%  " }
    26 = @{ D = "13.07";      E = "  +6.97%  " }
    27 = @{ D = "10.08";      E = "  +0.86%  " }
    28 = @{ D = $null;        E = "  -0.03%  " }
    29 = @{ D = "7.99";       E = "  -0.08%  " }
    30 = @{ D = "2.39";       E = "  -1.81%  " }
    31 = @{ D = $null;        E = "  +4.27%  " }
    34 = @{ D = $null;        E = "  +2.46%  " }
    35 = @{ D = "1.00";       E = $null }
    36 = @{ D = "1.01";       E = "  +3.78%  " }
    37 = @{ D = "5.92";       E = "  +3.43%  " }
    38 = @{ D = "48.47";      E = "  +7.13%  " }
    39 = @{ D = $null;        E = "  +8.27%  " }
    40 = @{ D = "50.24";      E = "  +2.30%  " }
    41 = @{ D = "0.318";      E = "  +5.55%  " }
    42 = @{ D = $null;        E = "  +0.43%  " }
    43 = @{ D = "8.69";       E = "  +1.96%  " }
    44 = @{ D = "2.80";       E = "  -1.10%  " }
    47 = @{ D = "381.39";     E = "  -0.77%  " }
    48 = @{ D = "135.49";     E = "  +1.69%  " }
    49 = @{ D = $null;        E = "  +0.00%  " }
    50 = @{ D = "25.03";      E = "  +5.94%  " }
    51 = @{ D = "2.24";       E = "  +3.50%  " }
}

foreach ($row in $updates.Keys) {
    $u = $updates[$row]
    if ($null -ne $u.D) {
        $dCell = $ws.Cells.Item($row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($row, 5).Value = $u.E
    }
}

# --- Rows that swapped coin identity/rank between snapshots ---------------
# Row 32/33: PEPE and EthereumClassic swapped order
$ws.Cells.Item(32, 2).Value = "EthereumClassic"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$d32 = $ws.Cells.Item(32, 4)
$d32.NumberFormat = "@"
$d32.Value = "29.03"
$ws.Cells.Item(32, 5).Value = "  +7.06%  "

$ws.Cells.Item(33, 2).Value = "PEPE"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$d33 = $ws.Cells.Item(33, 4)
$d33.NumberFormat = "@"
$d33.Value = "0.0000101"
$ws.Cells.Item(33, 5).Value = "  -0.23%  "

# Row 45/46: Maker and VeChain swapped order
$ws.Cells.Item(45, 2).Value = "VeChain"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$d45 = $ws.Cells.Item(45, 4)
$d45.NumberFormat = "@"
$d45.Value = "0.0362"
$ws.Cells.Item(45, 5).Value = "  +3.20%  "

$ws.Cells.Item(46, 2).Value = "Maker"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$d46 = $ws.Cells.Item(46, 4)
$d46.NumberFormat = "@"
$d46.Value = "2.843.00"
$ws.Cells.Item(46, 5).Value = "  +6.17%  "
